$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F25").Value = 101
$ws.Range("G25").Value = 2586.61
$ws.Range("F27").Value = 46
$ws.Range("G27").Value = 1885.54
$ws.Range("F31").Value = 88
$ws.Range("G31").Value = 3155.68
$ws.Range("B32").Value = 51755
$ws.Range("F32").Value = 1
$ws.Range("G32").Value = 30.74
$ws.Range("B33").Value = 66452
$ws.Range("F33").Value = 64
$ws.Range("G33").Value = 1967.36
$ws.Range("B40").Value = 71904.49000000001
$ws.Range("F102").Value = 2
$ws.Range("G102").Value = 4640.34
$ws.Range("B104").Value = 32658.18
$ws.Range("F119").Value = 31
$ws.Range("G119").Value = 2177.44
$ws.Range("B145").Value = 89412.44
$ws.Range("B148").Value = 65258
$ws.Range("F148").Value = 2
$ws.Range("G148").Value = 64287.16
$ws.Range("B149").Value = 64196
$ws.Range("F149").Value = 1
$ws.Range("G149").Value = 32143.58
$ws.Range("F221").Value = 80
$ws.Range("G221").Value = 6288.8
$ws.Range("F228").Value = 11
$ws.Range("G228").Value = 375.21
$ws.Range("F242").Value = 14
$ws.Range("G242").Value = 1246.84
$ws.Range("B247").Value = 89368.23
$ws.Range("B322").Value = 66188
$ws.Range("C322").Value = "HIM-Baby Care Gift Pack (Ww)1"
$ws.Range("D322").Value = 315.8
$ws.Range("E322").Value = 377.31
$ws.Range("F322").Value = 35
$ws.Range("G322").Value = 11053
$ws.Range("B323").Value = 48719
$ws.Range("C323").Value = "HIM-BABY CARE GIFT PACK (WW)1"
$ws.Range("D323").Value = 295.75
$ws.Range("E323").Value = 353.35
$ws.Range("F323").Value = -82
$ws.Range("G323").Value = -24251.5
$ws.Range("F349").Value = 10
$ws.Range("G349").Value = 759.1
$ws.Range("B367").Value = 64983
$ws.Range("C367").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F367").Value = 6
$ws.Range("G367").Value = 514.08
$ws.Range("B368").Value = 66194
$ws.Range("C368").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F368").Value = 35
$ws.Range("G368").Value = 2998.8
$ws.Range("B369").Value = 66196
$ws.Range("C369").Value = "HIM-Total Care Baby Pants Drapers-Xl-9S"
$ws.Range("F369").Value = 28
$ws.Range("G369").Value = 2455.6
$ws.Range("B370").Value = 64985
$ws.Range("C370").Value = "HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S"
$ws.Range("F370").Value = 13
$ws.Range("G370").Value = 1140.1
$ws.Range("B372").Value = 142347.22
$ws.Range("B397").Value = 60325
$ws.Range("E397").Value = 151.57
$ws.Range("F397").Value = -102
$ws.Range("G397").Value = -12939.72
$ws.Range("B398").Value = 63560
$ws.Range("E398").Value = 134.87
$ws.Range("F398").Value = 1
$ws.Range("G398").Value = 126.86
$ws.Range("F402").Value = 24
$ws.Range("G402").Value = 3307.68
$ws.Range("B409").Value = 32956.25
$ws.Range("B548").Value = 53602
$ws.Range("E548").Value = 15.69
$ws.Range("F548").Value = -232
$ws.Range("G548").Value = -3050.8
$ws.Range("B549").Value = 65068
$ws.Range("E549").Value = 13.97
$ws.Range("F549").Value = 0
$ws.Range("G549").Value = 0
$ws.Range("B556").Value = 45706
$ws.Range("E556").Value = 23.58
$ws.Range("F556").Value = -207
$ws.Range("G556").Value = -4084.11
$ws.Range("B557").Value = 64922
$ws.Range("E557").Value = 20.98
$ws.Range("F557").Value = 0
$ws.Range("G557").Value = 0
$ws.Range("B559").Value = 45718
$ws.Range("E559").Value = 19.38
$ws.Range("F559").Value = -295
$ws.Range("G559").Value = -4784.9
$ws.Range("B560").Value = 64927
$ws.Range("E560").Value = 17.26
$ws.Range("F560").Value = 0
$ws.Range("G560").Value = 0
$ws.Range("B569").Value = 65067
$ws.Range("E569").Value = 15.65
$ws.Range("F569").Value = 0
$ws.Range("G569").Value = 0
$ws.Range("B570").Value = 53595
$ws.Range("E570").Value = 17.61
$ws.Range("F570").Value = -338
$ws.Range("G570").Value = -4978.74
$ws.Range("F573").Value = 5
$ws.Range("G573").Value = 253.35
$ws.Range("B584").Value = 36161.87
$ws.Range("F608").Value = 111
$ws.Range("G608").Value = 31408.56
$ws.Range("F609").Value = 70
$ws.Range("G609").Value = 15556.1
$ws.Range("B612").Value = 127754.24
$ws.Range("B640").Value = 64810
$ws.Range("E640").Value = 291.22
$ws.Range("F640").Value = 2
$ws.Range("G640").Value = 547.84
$ws.Range("B641").Value = 53319
$ws.Range("E641").Value = 310.64
$ws.Range("F641").Value = -6
$ws.Range("G641").Value = -1643.52
$ws.Range("B669").Value = 64830
$ws.Range("E669").Value = 34.9
$ws.Range("F669").Value = 89
$ws.Range("G669").Value = 2921.87
$ws.Range("B670").Value = 60022
$ws.Range("E670").Value = 37.22
$ws.Range("F670").Value = -113
$ws.Range("G670").Value = -3709.79
$ws.Range("F685").Value = 503
$ws.Range("G685").Value = 27604.64
$ws.Range("F689").Value = 196
$ws.Range("G689").Value = 16779.56
$ws.Range("B692").Value = 160375.39
$ws.Range("F808").Value = 92
$ws.Range("G808").Value = 10010.52
$ws.Range("F832").Value = 390
$ws.Range("G832").Value = 14363.7
$ws.Range("B839").Value = 278195.82
$ws.Range("F878").Value = 85
$ws.Range("G878").Value = 6826.35
$ws.Range("F883").Value = 5
$ws.Range("G883").Value = 446.15
$ws.Range("B884").Value = 20208.25
$ws.Range("F890").Value = 1490
$ws.Range("G890").Value = 243033.9
$ws.Range("B896").Value = 269539.71
$ws.Range("F908").Value = 23
$ws.Range("G908").Value = 3660.68
$ws.Range("B912").Value = 16462.14
$ws.Range("F914").Value = 10
$ws.Range("G914").Value = 2423.9
$ws.Range("B935").Value = 91769.64
$ws.Range("B941").Value = 3949298.35
$ws.Range("B942").Value = 3949298.35
